# Hortaliza, Femacal de La Calera - Sandia
# A new weekly price-report row is inserted at row 186 (pushing the
# previously-existing rows 186-258 down to 187-259). The new row carries a
# fresh "Primera" quality entry dated 44488 (2021-10-19) sourced from Perú.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at 186; everything at/after 186 shifts down
# by one (old row 186 -> 187, ..., old row 258 -> 259). Formatting (e.g. the
# date style on column D) carries down with the shifted rows automatically.
$ws.Rows("186:186").Insert()

# Populate the newly inserted row 186 with the new record.
$ws.Range("A186").Value = 3
$ws.Range("B186").Value = "Femacal de La Calera"
$ws.Range("C186").Value = "Coquimbo"
$ws.Range("D186").Value = 44488
$ws.Range("E186").Value = 5
$ws.Range("F186").Value = 100112028
$ws.Range("G186").Value = "Sandia"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 250
$ws.Range("K186").Value = 800
$ws.Range("L186").Value = 800
$ws.Range("M186").Value = 800
$ws.Range("N186").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O186").Value = "Perú"
$ws.Range("P186").Value = 800
$ws.Range("Q186").Value = 1
$ws.Range("R186").Value = "Hortaliza"

# Make sure the new row's date cell keeps the same number format the other
# date cells in column D use (style index 2 / yyyy-mm-dd-ish date format).
$ws.Range("D186").NumberFormat = $ws.Range("D187").NumberFormat
